$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.438.11"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").Value = "3.448.65"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "407.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.596"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.50%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.691"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.129"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.04"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.73%  "
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("D15").Value = "3.463.06"
$ws.Range("E15").Value = "  -1.25%  "
$ws.Range("D16").Value = "62.515.60"
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.45%  "
$ws.Range("E18").Value = "  -2.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000145"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "84.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "313.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.87%  "
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.174"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.36%  "
$ws.Range("E31").Value = "  -3.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "43.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.58%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0488"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.997"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.322"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.50%  "
$ws.Range("E40").Value = "  -4.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "138.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.76%  "
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.03%  "
$ws.Range("D48").Value = "2.128.00"
$ws.Range("E48").Value = "  -3.44%  "
$ws.Range("E49").Value = "  -3.30%  "
$ws.Range("B50").Value = "Fetch.AI"
$ws.Range("C50").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +19.49%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.86%  "
